$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated in each data row: B, C, D, E, F, I, J, K, L, M (col indices 2,3,4,5,6,9,10,11,12,13)
$colIdx = @(2,3,4,5,6,9,10,11,12,13)

# New vm_pu values per row (rows 2-25), in column order matching $colIdx above
$rowsData = @{
    2 = @(1.02, 1.024525819036797, 1.027687710516158, 1.024933552443948, 1.023009657280859, 1.028540493957157, 1.029699942307616, 1.03050671189918, 1.027760601293353, 1.025842355591445)
    3 = @(1.02, 1.026161437136498, 1.028856640750343, 1.026344343702474, 1.025292674300836, 1.028848906323447, 1.030970959808858, 1.031482444217269, 1.028976949128002, 1.027928137798961)
    4 = @(1.02, 1.027216559348742, 1.029610172946366, 1.027254609729758, 1.026766015625233, 1.02904592541324, 1.03178989909272, 1.032110404993197, 1.029760893682251, 1.029273558623409)
    5 = @(1.02, 1.027659374512954, 1.029926286923337, 1.027636672956733, 1.027384496626049, 1.029128146644078, 1.03213335714089, 1.032373594298924, 1.03008972997923, 1.029838188597965)
    6 = @(1.02, 1.027733681114168, 1.029979324790005, 1.027700787484404, 1.02748828964073, 1.029141916539691, 1.032190977296248, 1.032417737981618, 1.030144900329333, 1.02993293542691)
    7 = @(1.02, 1.027222479224694, 1.02961439949967, 1.027259717267628, 1.026774283343722, 1.02904702643206, 1.031794491616827, 1.032113924893461, 1.029765290477369, 1.029281107072734)
    8 = @(1.02, 1.025079262460586, 1.02808335076466, 1.025410883861994, 1.023782042605536, 1.028645252316325, 1.030130219620264, 1.030837175534647, 1.028172323415551, 1.026548144470262)
    9 = @(1.02, 1.021277184534424, 1.025363220781516, 1.022132456322145, 1.018478039050548, 1.0279176286784, 1.027170236782255, 1.028560892051068, 1.02534094549041, 1.021698855559106)
    10 = @(1.02, 1.018724362017634, 1.023534227281019, 1.019932240805937, 1.01491924672758, 1.027419111115213, 1.025177751421224, 1.027024963942713, 1.023436247308853, 1.018441894768158)
    11 = @(1.02, 1.017614445821065, 1.022738417258138, 1.018975885074752, 1.013372452707878, 1.027200009760929, 1.024310259714223, 1.026355390141623, 1.022607268868115, 1.017025518474922)
    12 = @(1.02, 1.017201474107578, 1.022442228136365, 1.018620088630881, 1.012796997185136, 1.027118134505714, 1.023987308178796, 1.026105991807054, 1.022298699342182, 1.016498467029757)
    13 = @(1.02, 1.017290089904134, 1.022505788600579, 1.018696433916071, 1.012920475957933, 1.027135719329091, 1.024056615505821, 1.026159519956526, 1.022364918199222, 1.016611564677435)
    14 = @(1.02, 1.017580323830693, 1.022713946277783, 1.018946486410964, 1.013324904082716, 1.0271932519745, 1.024283579350206, 1.026334788917443, 1.022581775758439, 1.016981971722612)
    15 = @(1.02, 1.017759053271298, 1.022842120559513, 1.019100476896283, 1.013573964341206, 1.027228634504375, 1.024423322530948, 1.026442686281114, 1.022715302292909, 1.017210065404292)
    16 = @(1.02, 1.018797926479437, 1.023586960536372, 1.019995632829095, 1.015021776983073, 1.027433583506643, 1.025235222964499, 1.027069305401476, 1.023491173509505, 1.018535764023185)
    17 = @(1.02, 1.019448358622434, 1.024053141217993, 1.020556153124788, 1.015928371361886, 1.027561272023931, 1.025743228784279, 1.027461152205052, 1.023976714099404, 1.019365687213086)
    18 = @(1.02, 1.019827309191882, 1.024324686499199, 1.020882744204919, 1.016456613926625, 1.027635438273585, 1.026039084434718, 1.027689275776093, 1.02425951440155, 1.019849181822854)
    19 = @(1.02, 1.019956448254305, 1.024417214031584, 1.020994044259679, 1.016636637262416, 1.027660674210939, 1.026139886762494, 1.027766986749718, 1.024355873338586, 1.020013942607701)
    20 = @(1.02, 1.019378618593231, 1.024003162805108, 1.020496051009811, 1.015831160387785, 1.027547604589697, 1.02568877179887, 1.027419155730671, 1.023924662422042, 1.01927670509919)
    21 = @(1.02, 1.017494876636564, 1.022652665412852, 1.018872867858495, 1.013205835402476, 1.027176323641874, 1.02421676433522, 1.026283195686443, 1.022517934671761, 1.016872922452826)
    22 = @(1.02, 1.016306436949061, 1.021800135179709, 1.01784904112612, 1.011549922551789, 1.026940040039956, 1.023287043731903, 1.02556498133815, 1.021629701051387, 1.015356078569646)
    23 = @(1.02, 1.016936842306824, 1.02225240583533, 1.018392105887928, 1.012428263980413, 1.027065569593559, 1.023780310609135, 1.025946102627334, 1.022100932486619, 1.016160717149333)
    24 = @(1.02, 1.019410132455771, 1.024025747041799, 1.020523209652602, 1.015875087567518, 1.027553781279861, 1.025713379973826, 1.027438133470701, 1.023948183592812, 1.019316914094581)
    25 = @(1.02, 1.022263230495427, 1.026069137816638, 1.022982523262792, 1.019853138259466, 1.028108088081285, 1.027938782195218, 1.029152565777022, 1.026075886678193, 1.02295664205883)
}

foreach ($rowNum in $rowsData.Keys) {
    $rowData = $rowsData[$rowNum]
    for ($i = 0; $i -lt $colIdx.Length; $i++) {
        $ws.Cells.Item($rowNum, $colIdx[$i]).Value = $rowData[$i]
    }
}
